$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(28, 8).Value = 530.55554
$ws.Cells.Item(28, 9).Value = 590
$ws.Cells.Item(28, 10).Value = 233.33333
$ws.Cells.Item(28, 11).Value = 590
$ws.Cells.Item(28, 12).Value = 233.33333
$ws.Cells.Item(28, 13).Value = -105
$ws.Cells.Item(28, 14).Value = -1203.33333
$ws.Cells.Item(41, 8).Value = 312.95456
$ws.Cells.Item(41, 9).Value = 338.46155
$ws.Cells.Item(41, 10).Value = 276.1111
$ws.Cells.Item(41, 11).Value = 338.46155
$ws.Cells.Item(41, 12).Value = 276.1111
$ws.Cells.Item(41, 13).Value = 101.53845
$ws.Cells.Item(41, 14).Value = -1156.1111
$ws.Cells.Item(64, 8).Value = 3460.5945
$ws.Cells.Item(64, 9).Value = 3222.4546
$ws.Cells.Item(64, 10).Value = 3561.3462
$ws.Cells.Item(64, 11).Value = 3222.4546
$ws.Cells.Item(64, 12).Value = 3561.3462
$ws.Cells.Item(64, 13).Value = -2974.4546
$ws.Cells.Item(64, 14).Value = -4057.3462
$ws.Cells.Item(67, 8).Value = 3460.5945
$ws.Cells.Item(67, 9).Value = 3222.4546
$ws.Cells.Item(67, 10).Value = 3561.3462
$ws.Cells.Item(67, 11).Value = 3222.4546
$ws.Cells.Item(67, 12).Value = 3561.3462
$ws.Cells.Item(67, 13).Value = -2364.4546
$ws.Cells.Item(67, 14).Value = -5277.3462
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 763.60974
$ws.Cells.Item(2, 9).Value = 445.1613
$ws.Cells.Item(2, 11).Value = 445.1613
$ws.Cells.Item(2, 13).Value = -332.1613
$ws.Cells.Item(63, 8).Value = 3589.9443
$ws.Cells.Item(63, 9).Value = 3212.8823
$ws.Cells.Item(63, 10).Value = 10000
$ws.Cells.Item(63, 11).Value = 3212.8823
$ws.Cells.Item(63, 12).Value = 10000
$ws.Cells.Item(63, 13).Value = -2526.8823
$ws.Cells.Item(63, 14).Value = -11372
$ws.Cells.Item(66, 8).Value = 3589.9443
$ws.Cells.Item(66, 9).Value = 3212.8823
$ws.Cells.Item(66, 10).Value = 10000
$ws.Cells.Item(66, 11).Value = 16064.4115
$ws.Cells.Item(66, 12).Value = 50000
$ws.Cells.Item(66, 13).Value = -12632.4115
$ws.Cells.Item(66, 14).Value = -56864
$ws.Cells.Item(110, 8).Value = 1080.8
$ws.Cells.Item(110, 9).Value = 1127.421
$ws.Cells.Item(110, 10).Value = 933.1667
$ws.Cells.Item(110, 11).Value = 1127.421
$ws.Cells.Item(110, 12).Value = 933.1667
$ws.Cells.Item(110, 13).Value = 917.579
$ws.Cells.Item(110, 14).Value = -5023.1667
$ws.Cells.Item(116, 8).Value = 763.60974
$ws.Cells.Item(116, 9).Value = 445.1613
$ws.Cells.Item(116, 11).Value = 445.1613
$ws.Cells.Item(116, 13).Value = 1848.8387
$ws.Cells.Item(124, 8).Value = 23532.666
$ws.Cells.Item(124, 10).Value = 23532.666
$ws.Cells.Item(124, 12).Value = 23532.666
$ws.Cells.Item(124, 14).Value = -33352.666
$ws.Cells.Item(125, 8).Value = 72415.836
$ws.Cells.Item(125, 10).Value = 72415.836
$ws.Cells.Item(125, 12).Value = 72415.836
$ws.Cells.Item(125, 14).Value = -82255.836
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 763.60974
$ws.Cells.Item(3, 9).Value = 445.1613
$ws.Cells.Item(3, 11).Value = 445.1613
$ws.Cells.Item(3, 13).Value = -331.1613
$ws.Cells.Item(105, 8).Value = 2166.5
$ws.Cells.Item(105, 9).Value = 2166.5
$ws.Cells.Item(105, 11).Value = 2166.5
$ws.Cells.Item(105, 13).Value = -419.5
$ws.Cells.Item(107, 8).Value = 26631.182
$ws.Cells.Item(107, 9).Value = 27699.334
$ws.Cells.Item(107, 11).Value = 27699.334
$ws.Cells.Item(107, 13).Value = -25779.334
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 1297.4286
$ws.Cells.Item(16, 9).Value = 1180.25
$ws.Cells.Item(16, 10).Value = 2000.5
$ws.Cells.Item(16, 11).Value = 1180.25
$ws.Cells.Item(16, 12).Value = 2000.5
$ws.Cells.Item(16, 13).Value = -893.25
$ws.Cells.Item(16, 14).Value = -2574.5
$ws.Cells.Item(107, 8).Value = 465.77777
$ws.Cells.Item(107, 9).Value = 432.3125
$ws.Cells.Item(107, 10).Value = 514.4545000000001
$ws.Cells.Item(107, 11).Value = 432.3125
$ws.Cells.Item(107, 12).Value = 514.4545000000001
$ws.Cells.Item(107, 13).Value = 1487.6875
$ws.Cells.Item(107, 14).Value = -4354.4545
$ws.Cells.Item(113, 8).Value = 1297.4286
$ws.Cells.Item(113, 9).Value = 1180.25
$ws.Cells.Item(113, 10).Value = 2000.5
$ws.Cells.Item(113, 11).Value = 1180.25
$ws.Cells.Item(113, 12).Value = 2000.5
$ws.Cells.Item(113, 13).Value = 989.75
$ws.Cells.Item(113, 14).Value = -6340.5
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(131, 8).Value = 868.79
$ws.Cells.Item(131, 10).Value = 884.57294
$ws.Cells.Item(131, 12).Value = 2653.71882
$ws.Cells.Item(131, 14).Value = -12733.71882
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 5881.1816
$ws.Cells.Item(70, 9).Value = 5634.364
$ws.Cells.Item(70, 10).Value = 6128
$ws.Cells.Item(70, 11).Value = 5634.364
$ws.Cells.Item(70, 12).Value = 6128
$ws.Cells.Item(70, 13).Value = -5364.364
$ws.Cells.Item(70, 14).Value = -6668
$ws.Cells.Item(73, 8).Value = 5881.1816
$ws.Cells.Item(73, 9).Value = 5634.364
$ws.Cells.Item(73, 10).Value = 6128
$ws.Cells.Item(73, 11).Value = 5634.364
$ws.Cells.Item(73, 12).Value = 6128
$ws.Cells.Item(73, 13).Value = -4698.364
$ws.Cells.Item(73, 14).Value = -8000
$ws.Cells.Item(80, 8).Value = 3346.4614
$ws.Cells.Item(80, 9).Value = 3625.625
$ws.Cells.Item(80, 10).Value = 2899.8
$ws.Cells.Item(80, 11).Value = 3625.625
$ws.Cells.Item(80, 12).Value = 2899.8
$ws.Cells.Item(80, 13).Value = -2627.625
$ws.Cells.Item(80, 14).Value = -4895.8
$ws.Cells.Item(83, 8).Value = 3346.4614
$ws.Cells.Item(83, 9).Value = 3625.625
$ws.Cells.Item(83, 10).Value = 2899.8
$ws.Cells.Item(83, 11).Value = 18128.125
$ws.Cells.Item(83, 12).Value = 14499
$ws.Cells.Item(83, 13).Value = -13136.125
$ws.Cells.Item(83, 14).Value = -24483
$ws.Cells.Item(102, 8).Value = 1692.881
$ws.Cells.Item(102, 9).Value = 1677.8857
$ws.Cells.Item(102, 10).Value = 1767.8572
$ws.Cells.Item(102, 11).Value = 1677.8857
$ws.Cells.Item(102, 12).Value = 1767.8572
$ws.Cells.Item(102, 13).Value = -55.88570000000004
$ws.Cells.Item(102, 14).Value = -5011.8572
$ws.Cells.Item(113, 8).Value = 1926.8462
$ws.Cells.Item(113, 9).Value = 2251.25
$ws.Cells.Item(113, 10).Value = 1407.8
$ws.Cells.Item(113, 11).Value = 2251.25
$ws.Cells.Item(113, 12).Value = 1407.8
$ws.Cells.Item(113, 13).Value = -81.25
$ws.Cells.Item(113, 14).Value = -5747.8
$ws.Cells.Item(126, 8).Value = 2412.125
$ws.Cells.Item(126, 9).Value = 2252.8125
$ws.Cells.Item(126, 10).Value = 2730.75
$ws.Cells.Item(126, 11).Value = 6758.4375
$ws.Cells.Item(126, 12).Value = 8192.25
$ws.Cells.Item(126, 13).Value = -4288.4375
$ws.Cells.Item(126, 14).Value = -13132.25
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(61, 8).Value = 18569.54
$ws.Cells.Item(61, 9).Value = 25500.555
$ws.Cells.Item(61, 10).Value = 2974.75
$ws.Cells.Item(61, 11).Value = 25500.555
$ws.Cells.Item(61, 12).Value = 2974.75
$ws.Cells.Item(61, 13).Value = -25298.555
$ws.Cells.Item(61, 14).Value = -3378.75
$ws.Cells.Item(113, 8).Value = 18569.54
$ws.Cells.Item(113, 9).Value = 25500.555
$ws.Cells.Item(113, 10).Value = 2974.75
$ws.Cells.Item(113, 11).Value = 25500.555
$ws.Cells.Item(113, 12).Value = 2974.75
$ws.Cells.Item(113, 13).Value = -23330.555
$ws.Cells.Item(113, 14).Value = -7314.75
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(100, 8).Value = 21069.2
$ws.Cells.Item(100, 9).Value = 33868.168
$ws.Cells.Item(100, 11).Value = 67736.336
$ws.Cells.Item(100, 13).Value = -67195.336
$ws.Cells.Item(107, 8).Value = 430.92856
$ws.Cells.Item(107, 9).Value = 288.3158
$ws.Cells.Item(107, 10).Value = 732
$ws.Cells.Item(107, 11).Value = 864.9474
$ws.Cells.Item(107, 12).Value = 2196
$ws.Cells.Item(107, 13).Value = 1055.0526
$ws.Cells.Item(107, 14).Value = -6036
$ws.Cells.Item(113, 8).Value = 748.25806
$ws.Cells.Item(113, 9).Value = 486.26086
$ws.Cells.Item(113, 10).Value = 1501.5
$ws.Cells.Item(113, 11).Value = 1458.78258
$ws.Cells.Item(113, 12).Value = 4504.5
$ws.Cells.Item(113, 13).Value = 711.2174199999999
$ws.Cells.Item(113, 14).Value = -8844.5
